$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "motivos_venta"

# Reset header cell values (lowercase, no accents) and clear any formatting/styles
$ws.Range("A1").Value = "codigo"
$ws.Range("B1").Value = "descripcion"
$ws.Range("C1").Value = "comentario"

$ws.Range("A1:C1").ClearFormats()
$ws.Columns("A:C").ColumnWidth = 8.43
